# Pagination_Generator.xlsx update:
#  - Add new HTML pages for Section 3 (Structure), Section 4 (Links) and
#    Section 6 (Images) to the Table1 data table on the Pagination_Gen sheet.
#  - Re-point the "Current file" dropdown (C2) at the new last page
#    (6images_5text.html) which now drives the prev/next/current output.
#  - Re-apply the Section filter so it shows the new current section (6)
#    instead of the old one (5).
#  - Grow the "filename" named range and the Table1 ListObject to cover the
#    newly added rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pagination_Gen")
$tbl = $ws.ListObjects.Item("Table1")

$sectionPageFormula = '=[@Section]&"_"&[@[Page #2]]'
$outputFormula = '=IF([@[File Name]]=current_file,current_left&[@[File Name]]&current_mid&[@[Page ''#2]]&current_right,other_left&[@[File Name]]&other_mid&[@[Page ''#2]]&other_right)'

# ---------------------------------------------------------------------------
# 1. Make room in the worksheet for the new rows.
#    Layout before:  row10 header | row11 previous | rows12-26 sec1/2 |
#                     rows27-34 sec5 | row35 next
#    Layout after:   row10 header | row11 previous | rows12-26 sec1/2 |
#                     rows27-30 sec3 | rows31-32 sec4 | rows33-40 sec5 |
#                     rows41-45 sec6 | row46 next
# ---------------------------------------------------------------------------

# Room for Section 3 (4 rows) + Section 4 (2 rows) right before the old
# Section 5 block (worksheet rows 27-32 today).
$ws.Range("B27:F32").Insert()

# Room for Section 6 (5 rows) right before the trailing "next" row, which
# (after the insert above) now sits at worksheet row 41.
$ws.Range("B41:F45").Insert()

# Grow the table to cover everything through the new "next" row (46).
$tbl.Resize($ws.Range("B10:F46"))

# ---------------------------------------------------------------------------
# 2. Fill in the new data rows.
# ---------------------------------------------------------------------------

function Fill-Row($rowNum, $file, $section, $page) {
    $ws.Range("C$rowNum").Value = $file
    $ws.Range("D$rowNum").Value = $section
    $ws.Range("E$rowNum").Value = $page
    $ws.Range("B$rowNum").Formula = $sectionPageFormula
    $ws.Range("F$rowNum").Formula = $outputFormula
}

# Section 3 - Structure (hidden rows, filter stays on the current section)
Fill-Row 27 "3structure_1heading.html" 3 1
Fill-Row 28 "3structure_2heading-det.html" 3 2
Fill-Row 29 "3structure_3heading-lev.html" 3 3
Fill-Row 30 "3structure_4list.html" 3 4

# Section 4 - Links
Fill-Row 31 "4links_1link-purpose.html" 4 1
Fill-Row 32 "4links_2change.html" 4 2

# Section 6 - Images (these become the visible rows once the filter flips
# over to section 6 below)
Fill-Row 41 "6images_1meaningful.html" 6 1
Fill-Row 42 "6images_2decorative.html" 6 2
Fill-Row 43 "6images_3background.html" 6 3
Fill-Row 44 "6images_4captcha.html" 6 4
Fill-Row 45 "6images_5text.html" 6 5

# ---------------------------------------------------------------------------
# 3. Point "Current file" at the new final page and let the Section/Page
#    lookups (D2/E2) recompute from it.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "6images_5text.html"

# ---------------------------------------------------------------------------
# 4. Grow the "filename" named range (dropdown source for C2) to cover the
#    new rows.
# ---------------------------------------------------------------------------
$wb.Names.Item("filename").RefersTo = "=Pagination_Gen!`$C`$11:`$C`$46"

# ---------------------------------------------------------------------------
# 5. Re-apply the table's Section filter for the new current section (6),
#    keeping the "prev-next" helper rows visible too.
# ---------------------------------------------------------------------------
$tbl.Range.AutoFilter(3, @("6", "prev-next"), 7)

# ---------------------------------------------------------------------------
# 6. Restore the active selection used when the sheet was last saved.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("F11:F45").Select()

$excel.CalculateFull()
